$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Games")
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 1
$win.ScrollRow = 104
Write-Host ("SplitRow=" + $win.SplitRow())
Write-Host ("ScrollRow=" + $win.ScrollRow())
$ws.Range("AU128").Select()
